# Refresh the cryptos.xlsx price/volume/link snapshot (GitHub Actions data pull).
# Source: commit "Updated cryptos list on Sun Apr 21 20:27:05 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '64.615.68'
$ws.Range('E2').Value = '  -0.14%  '

# Row 3
$ws.Range('D3').Value = '3.145.59'
$ws.Range('E3').Value = '  -0.13%  '

# Row 4
$ws.Range('E4').Value = '  -0.05%  '

# Row 5
$ws.Range('D5').Value = "'576.84"
$ws.Range('E5').Value = '  +0.83%  '

# Row 6
$ws.Range('D6').Value = "'148.69"
$ws.Range('E6').Value = '  -1.42%  '

# Row 7
$ws.Range('D7').Value = "'0.999"
$ws.Range('E7').Value = '  -0.14%  '

# Row 8
$ws.Range('D8').Value = '3.146.43'
$ws.Range('E8').Value = '  +0.04%  '

# Row 9
$ws.Range('D9').Value = "'0.525"
$ws.Range('E9').Value = '  -0.18%  '

# Row 10
$ws.Range('E10').Value = '  -2.30%  '

# Row 11
$ws.Range('D11').Value = "'6.12"
$ws.Range('E11').Value = '  -0.80%  '

# Row 12
$ws.Range('D12').Value = "'0.501"
$ws.Range('E12').Value = '  -0.27%  '

# Row 13
$ws.Range('D13').Value = "'0.0000261"
$ws.Range('E13').Value = '  +2.17%  '

# Row 14
$ws.Range('D14').Value = "'37.11"
$ws.Range('E14').Value = '  -1.57%  '

# Row 15
$ws.Range('D15').Value = '3.665.62'
$ws.Range('E15').Value = '  -0.22%  '

# Row 16
$ws.Range('D16').Value = '64.711.50'
$ws.Range('E16').Value = '  -0.23%  '

# Row 17
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').Value = "'7.12"
$ws.Range('E17').Value = '  -1.32%  '

# Row 18
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.133.75'
$ws.Range('E18').Value = '  -0.72%  '

# Row 19
$ws.Range('E19').Value = '  +0.36%  '

# Row 20
$ws.Range('D20').Value = "'503.62"
$ws.Range('E20').Value = '  -1.53%  '

# Row 21
$ws.Range('D21').Value = "'14.85"
$ws.Range('E21').Value = '  -0.57%  '

# Row 22
$ws.Range('B22').Value = 'Polygon'
$ws.Range('C22').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D22').Value = "'0.712"
$ws.Range('E22').Value = '  -2.94%  '

# Row 23
$ws.Range('B23').Value = 'InternetComputer(DFINITY)'
$ws.Range('C23').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D23').Value = "'15.21"
$ws.Range('E23').Value = '  +0.36%  '

# Row 24
$ws.Range('D24').Value = "'7.70"
$ws.Range('E24').Value = '  -1.67%  '

# Row 25
$ws.Range('D25').Value = "'84.08"
$ws.Range('E25').Value = '  -1.20%  '

# Row 26
$ws.Range('E26').Value = '  +0.22%  '

# Row 27
$ws.Range('E27').Value = '  +2.01%  '

# Row 28
$ws.Range('D28').Value = "'2.89"
$ws.Range('E28').Value = '  -1.12%  '

# Row 29
$ws.Range('E29').Value = '  -0.75%  '

# Row 30
$ws.Range('E30').Value = '  +5.60%  '

# Row 31
$ws.Range('D31').Value = "'27.50"
$ws.Range('E31').Value = '  -1.70%  '

# Row 32
$ws.Range('D32').Value = "'0.999"
$ws.Range('E32').Value = '  -0.26%  '

# Row 33
$ws.Range('E33').Value = '  +0.61%  '

# Row 34
$ws.Range('D34').Value = "'6.17"
$ws.Range('E34').Value = '  +1.45%  '

# Row 35
$ws.Range('D35').Value = "'6.45"
$ws.Range('E35').Value = '  -1.62%  '

# Row 36
$ws.Range('D36').Value = "'54.57"

# Row 37
$ws.Range('B37').Value = 'Bittensor'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D37').Value = "'481.50"
$ws.Range('E37').Value = '  -0.23%  '

# Row 38
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = "'0.0889"
$ws.Range('E38').Value = '  +2.98%  '

# Row 39
$ws.Range('E39').Value = '  -1.87%  '

# Row 40
$ws.Range('E40').Value = '  -2.90%  '

# Row 41
$ws.Range('D41').Value = "'8.69"
$ws.Range('E41').Value = '  +0.73%  '

# Row 42
$ws.Range('D42').Value = '2.997.56'
$ws.Range('E42').Value = '  -3.69%  '

# Row 43
$ws.Range('E43').Value = '  -4.40%  '

# Row 44
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').Value = "'0.281"
$ws.Range('E44').Value = '  -4.40%  '

# Row 45
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').Value = "'2.41"
$ws.Range('E45').Value = '  -1.92%  '

# Row 46
$ws.Range('D46').Value = "'27.97"
$ws.Range('E46').Value = '  -4.01%  '

# Row 47
$ws.Range('D47').Value = '0.0₃0580'
$ws.Range('E47').Value = '  +0.33%  '

# Row 48
$ws.Range('E48').Value = '  -0.02%  '

# Row 49
$ws.Range('E49').Value = '  -1.52%  '

# Row 50
$ws.Range('E50').Value = '  -2.75%  '

# Row 51
$ws.Range('D51').Value = "'33.25"
$ws.Range('E51').Value = '  +5.26%  '
